$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.390"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05929"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.433"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8070"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9057"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1417"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07440"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03222"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03041"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09316"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.948"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001589"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04785"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006134"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").Value = "'0.007493"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004391"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009846"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00007814"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.609"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.150"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01104"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "24OneONE"
$ws.Range("D26").Value = "'0.3248"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.1322"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03864"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006210"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1060"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002805"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007262"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005186"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.0005809"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Value = "'0.002264"
$ws.Range("D49").Style = "Normal"
